# Trade #25 (MarketMaking strategy) closed at 2026-02-17 20:54:03 - unknown UNKNOWN +0.000%
# Updates the summary/status rollups, flips the open trade row to CLOSED with
# its exit figures, and appends the newly-opened trade row (#86) to both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: portfolio-level rollups
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.51   # Current Capital
$summary.Range("B4").Value = 0.31      # Total P&L $
$summary.Range("B5").Value = 0.12      # Total P&L %
$summary.Range("B6").Value = 53        # Total Trades
$summary.Range("B7").Value = 25        # Winning Trades
$summary.Range("B9").Value = 47.17     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.51     # Capital
$status.Range("D5").Value = 20         # Trades
$status.Range("E5").Value = 0.2        # P&L $
$status.Range("F5").Value = 0.51       # P&L %
$status.Range("G5").Value = 55         # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet: close out trade row 54 (Trade #53) and append new row 87
# (Trade #86). Column order: A Trade#, B Date, C Time, D Strategy, E Side,
# F Entry Price, G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
# L Exit Reason, M Duration (min), N Entry Slippage, O Exit Slippage,
# P Confidence, Q Entry Reason.
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G54").Value = 0.165017
$allTrades.Range("H54").Value = "CLOSED"
$allTrades.Range("I54").Value = 26.9358
$allTrades.Range("J54").Value = 0.04
$allTrades.Range("K54").Value = 100.51
$allTrades.Range("L54").Value = "early_exit"
$allTrades.Range("M54").Value = 0.13

$allTrades.Range("A87").Value = 86
# "2026-02-17" looks like a date to Excel's smart-entry parser; force it to
# stay literal text (matching the rest of the Date column) and strip the
# number-format trace it leaves behind.
$allTrades.Range("B87").NumberFormat = "@"
$allTrades.Range("B87").Value = "2026-02-17"
$allTrades.Range("B87").ClearFormats()
$allTrades.Range("C87").Value = "20:53:56"
$allTrades.Range("D87").Value = "MarketMaking"
$allTrades.Range("E87").Value = "UP"
$allTrades.Range("F87").Value = 0.13
$allTrades.Range("H87").Value = "OPEN"
$allTrades.Range("I87").Value = 0
$allTrades.Range("J87").Value = 0
$allTrades.Range("K87").Value = 100.4784370824165
$allTrades.Range("M87").Value = 0
$allTrades.Range("N87").Value = 0
$allTrades.Range("O87").Value = 0
$allTrades.Range("P87").Value = 0.6
$allTrades.Range("Q87").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet: close out trade row 21 (Trade #53) and append new row 54
# (Trade #86). Column order here differs from "All Trades": A Trade#, B Date,
# C Time, D Strategy, E Side, F Entry Price, G Exit Price, H Status,
# I P&L %, J P&L $, K Capital After, L Entry Slippage, M Exit Slippage,
# N Confidence, O Entry Reason, P Exit Reason, Q Duration (min).
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("G21").Value = 0.165017
$marketMaking.Range("H21").Value = "CLOSED"
$marketMaking.Range("I21").Value = 26.9358
$marketMaking.Range("J21").Value = 0.04
$marketMaking.Range("K21").Value = 100.51
$marketMaking.Range("P21").Value = "early_exit"
$marketMaking.Range("Q21").Value = 0.13

$marketMaking.Range("A54").Value = 86
$marketMaking.Range("B54").NumberFormat = "@"
$marketMaking.Range("B54").Value = "2026-02-17"
$marketMaking.Range("B54").ClearFormats()
$marketMaking.Range("C54").Value = "20:53:56"
$marketMaking.Range("D54").Value = "MarketMaking"
$marketMaking.Range("E54").Value = "UP"
$marketMaking.Range("F54").Value = 0.13
$marketMaking.Range("H54").Value = "OPEN"
$marketMaking.Range("I54").Value = 0
$marketMaking.Range("J54").Value = 0
$marketMaking.Range("K54").Value = 100.4784370824165
$marketMaking.Range("L54").Value = 0
$marketMaking.Range("M54").Value = 0
$marketMaking.Range("N54").Value = 0.6
$marketMaking.Range("O54").Value = "Normal spread capture: 19600 bps"
$marketMaking.Range("Q54").Value = 0
